# Apply updated parameter values and refresh selections across several sheets.

$wb = $excel.ActiveWorkbook

# --- Workspace: just move the active selection ---
$wsWorkspace = $wb.Worksheets.Item("Workspace")
$wsWorkspace.Activate()
$wsWorkspace.Range("E13").Select()

# --- 1.1_TumourAgnosticCosts: update raw costs (col C) and move selection ---
$ws1 = $wb.Worksheets.Item("1.1_TumourAgnosticCosts")
$ws1.Range("C3").Value = 139.41999999999999
$ws1.Range("C4").Value = 104.86
$ws1.Range("C5").Value = 133.49
$ws1.Range("C6").Value = 1494.01
$ws1.Range("C9").Value = 1494.01

# C8 / E8 lose their formulas and become plain literal values
$ws1.Range("C8").Value = 2002
$ws1.Range("E8").Value = 2002

$ws1.Activate()
$ws1.Range("A8:H8").Select()

# --- 1.3_TreatmentCost: just move the active selection ---
$ws2 = $wb.Worksheets.Item("1.3_TreatmentCost")
$ws2.Activate()
$ws2.Range("B8").Select()

# --- 1.4_AdminCost: update raw costs (col B) and move selection ---
$ws3 = $wb.Worksheets.Item("1.4_AdminCost")
$ws3.Range("B2").Value = 139.01
$ws3.Range("B3").Value = 116.94
$ws3.Range("B5").Value = 142.52000000000001
$ws3.Range("B6").Value = 183.3
$ws3.Range("B7").Value = 139.46
$ws3.Range("B9").Value = 3.07

$ws3.Activate()
$ws3.Range("B9").Select()

# --- 2.1_Utilities: just move the active selection ---
$ws4 = $wb.Worksheets.Item("2.1_Utilities")
$ws4.Activate()
$ws4.Range("B3:B4").Select()

# Leave the originally active sheet (Workspace) selected/active at the end,
# matching tabSelected="1" in the workbook.
$wsWorkspace.Activate()
$wsWorkspace.Range("E13").Select()
